$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 values
$ws.Range("A5").Value = "Arrays"
$ws.Range("B5").Value = 271
$ws.Range("C5").Value = "Encode and Decode Strings"
$ws.Range("D5").Value = "Encoding: Calculate a 4 string format for length of string. Append Length+string for each of the string oin list`nDecoding : First slice the 4 string to extract length of the string and then select the string. Move on to the next Length+string combo and repeat"
$ws.Range("E5").Value = "Time Complexity: O(n + L) for encoding and O(L) for decoding`nwhere`nn is the number of strings`nL is the total length of all strings combined"
$ws.Range("F5").Value = "Space Complexity: O(L + n) for encoding and O(L) for decoding"

# Match formatting of existing rows (row height + wrap text for whole row)
$ws.Range("A5:F5").RowHeight = 87
$ws.Range("A5:F5").WrapText = $true

# Hyperlink for C5, mirroring the other "Name" cells that link to LeetCode
$ws.Hyperlinks.Add($ws.Range("C5"), "https://leetcode.com/problems/encode-and-decode-strings", "", "", "https://leetcode.com/problems/encode-and-decode-strings")
$ws.Range("C5").HorizontalAlignment = -4131
$ws.Range("C5").VerticalAlignment = -4108
$ws.Range("C5").WrapText = $true

# Selection ends up parked past the data, as in the saved workbook
$ws.Range("E7").Select()
